$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20; existing rows 20-118 shift down to 21-119.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new data point.
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value = 45243
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 100112022
$ws.Range("G20").Value = "Arveja Verde"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 21000
$ws.Range("L20").Value = 21000
$ws.Range("M20").Value = 21000
$ws.Range("N20").Value = '$/malla 25 kilos'
$ws.Range("O20").Value = "Provincia de Limarí"
$ws.Range("P20").Value = 840
$ws.Range("Q20").Value = 25
$ws.Range("R20").Value = "Hortaliza"

# Match the date-number format used by the rest of column D.
$ws.Range("D20").NumberFormat = $ws.Range("D21").NumberFormat
